$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84; this shifts existing rows 84..167 down to 85..168
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly data point
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C84").Value = "Los Lagos"
$ws.Range("D84").Value = 45128
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = 100112022
$ws.Range("G84").Value = "Arveja Verde"
$ws.Range("H84").Value = "Perfection"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 45
$ws.Range("K84").Value = 40000
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = 40000
$ws.Range("N84").Value = "$/malla 25 kilos"
$ws.Range("O84").Value = "Provincia de Limarí"
$ws.Range("P84").Value = 1600
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
